$wb = $excel.ActiveWorkbook

# --- Sheet "ERT_SU_CZ": rename header label and drop the United Kingdom row ---
$ws1 = $wb.Worksheets.Item("ERT_SU_CZ")

# "SES Area" -> "SES Area (RP3)"
$ws1.Range("A6").Value = "SES Area (RP3)"

# United Kingdom was row 36 (last data row) - remove it entirely; the summary
# row 6 SUM()/shared formulas will auto-shrink their ranges to B7:B35 etc.
$ws1.Rows.Item(36).Delete()

# --- Sheet "Change Log": record the UK removal ---
$ws2 = $wb.Worksheets.Item("Change Log")

$ws2.Range("A2").Value = 44351
$ws2.Range("B2").Value = "UK"
$ws2.Range("C2").Value = 2020
$ws2.Range("D2").Value = "UK removed from RP3 area"
